# Generate Report for Handback
# - Overview!B3 / Overview!C3 shared text "Ready for handoff" -> "Handback transform failed"
# - zh-cn sheet, row 3: add L3 diagnostic message
# - de-de sheet, row 3: add L3 diagnostic message

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Handback transform failed"
$wsOverview.Range("C3").Value = "Handback transform failed"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Handback transform failed"
$wsZhCn.Range("L3").Value = "Handback file name: nfstk1t2.wgv is different with handoff file name: d3827935-abc4-4ab3-8e4c-1fbfc64c2448.4abfb4c43acccebca53155b4fbf32036aeecd0e7.zh-cn."

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Handback transform failed"
$wsDeDe.Range("L3").Value = "Handback file name: nfstk1t2.wgv is different with handoff file name: d3827935-abc4-4ab3-8e4c-1fbfc64c2448.4abfb4c43acccebca53155b4fbf32036aeecd0e7.de-de."
